$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: Kelly Marinduque / Referral / Fiber Technician
$ws.Cells.Item(9, 1).Value = 42920
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(9, 1).PasteSpecial(-4122)
$ws.Cells.Item(9, 2).Value = "Kelly Marinduque"
$ws.Cells.Item(9, 3).Value = "Referral"
$ws.Cells.Item(9, 4).Value = "Fiber Technician"
$ws.Cells.Item(9, 5).Value = 9233233212

# Row 10: Harry Potter / HOIT Website / System Administrator
$ws.Cells.Item(10, 1).Value = 42920
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(10, 1).PasteSpecial(-4122)
$ws.Cells.Item(10, 2).Value = "Harry Potter"
$ws.Cells.Item(10, 3).Value = "HOIT Website"
$ws.Cells.Item(10, 4).Value = "System Administrator"
$ws.Cells.Item(10, 5).Value = 9233233212

$excel.CutCopyMode = $false
